$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "0.999", "210.24") rather than
# numbers. Mark the cells whose new value would otherwise be auto-coerced
# to a number as Text first, so the literal string is preserved exactly,
# matching the inline-string cell type used by the source data feed.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '28.653.60'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.564.45'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '210.24'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").Value = '0.487'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '25.28'
$ws.Range("E8").Value = '  +6.63%  '
$ws.Range("D9").Value = '0.245'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").Value = '0.0585'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '1.788.57'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '1.561.55'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '28.670.53'
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '61.44'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '231.82'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("D19").Value = '7.37'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '0.0₃0675'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '0.996'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").Value = '3.91'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").Value = '9.01'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  +2.86%  '
$ws.Range("D25").Value = '150.74'
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("D26").Value = '14.80'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '6.22'
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("D30").Value = '0.0461'
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").Value = '3.16'
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("D33").Value = '1.390.71'
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  -3.93%  '
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.29'
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.65'
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("D41").Value = '0.517'
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("D43").Value = '0.776'
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").Value = '0.0459'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").Value = '64.10'
$ws.Range("E45").Value = '  +2.86%  '
$ws.Range("D46").Value = '5.28'
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("D47").Value = '1.702.34'
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("E48").Value = '  -5.60%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '85.26'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '43.35'
$ws.Range("E50").Value = '  +6.73%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0101'
$ws.Range("E51").Value = '  +0.04%  '
